$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B56 currently holds "3" stored as text; convert it to a real number (3),
# matching every other row's politeness_score column.
$ws.Range("B56").Value = 3

# Append a new annotation row (row 57) for Ying Tang.
$ws.Range("A57").Value = "Ying Tang"
# politeness_score is kept as text "3" (quote-prefixed) for this row, as in
# the source data, rather than converted to a number.
$ws.Range("B57").Value = "'3"
$ws.Range("B57").Style = "Normal"
$ws.Range("C57").Value = "无"
$ws.Range("D57").Value = "DIS"
$ws.Range("E57").Value = "MET"
$ws.Range("F57").Value = "a903e5ac-dd38-46eb-9fca-86d5e31ee0d2"
$ws.Range("G57").Value = "H1aIuk-RW_annotated.xlsx"
$ws.Range("H57").Value = "2) The ""Active learning"" approach is simply the classing hitting set approach for computing k-center."
